$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 71054.11
$ws.Range("J134").Value = 71054.11
$ws.Range("L134").Value = 71054.11
$ws.Range("N134").Value = -81194.11

$ws.Range("H138").Value = 3386.4736
$ws.Range("I138").Value = 5548.2
$ws.Range("J138").Value = 2614.4285
$ws.Range("K138").Value = 16644.6
$ws.Range("L138").Value = 7843.2855
$ws.Range("M138").Value = -11504.6
$ws.Range("N138").Value = -18123.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()

$ws.Range("H61").Value = 14616.193
$ws.Range("I61").Value = 3511.7
$ws.Range("K61").Value = 3511.7
$ws.Range("M61").Value = -3299.7

$ws.Range("H74").Value = 15516.069
$ws.Range("I74").Value = 826.625
$ws.Range("J74").Value = 33595.383
$ws.Range("K74").Value = 826.625
$ws.Range("L74").Value = 33595.383
$ws.Range("M74").Value = 47.375
$ws.Range("N74").Value = -35343.383

$ws.Range("H77").Value = 15516.069
$ws.Range("I77").Value = 826.625
$ws.Range("J77").Value = 33595.383
$ws.Range("K77").Value = 4133.125
$ws.Range("L77").Value = 167976.915
$ws.Range("M77").Value = 234.875
$ws.Range("N77").Value = -176712.915

$ws.Range("H108").Value = 49995
$ws.Range("J108").Value = 49995
$ws.Range("L108").Value = 49995
$ws.Range("N108").Value = -57675

$ws.Range("H132").Value = 3586998
$ws.Range("I132").Value = 5501.1333
$ws.Range("K132").Value = 16503.3999
$ws.Range("M132").Value = -13973.3999

$ws.Range("H136").Value = 14616.193
$ws.Range("I136").Value = 3511.7
$ws.Range("K136").Value = 10535.1
$ws.Range("M136").Value = -7985.099999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 25975.205
$ws.Range("I20").Value = 6673.7144
$ws.Range("K20").Value = 6673.7144
$ws.Range("M20").Value = -6426.7144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 18995.25
$ws.Range("J18").Value = 18995.25
$ws.Range("L18").Value = 18995.25
$ws.Range("N18").Value = -19455.25

$ws.Range("H58").Value = 18773.84
$ws.Range("I58").Value = 8648.154
$ws.Range("K58").Value = 8648.154
$ws.Range("M58").Value = -8445.154

$ws.Range("H99").Value = 9425.48
$ws.Range("I99").Value = 3540
$ws.Range("J99").Value = 10896.85
$ws.Range("K99").Value = 3540
$ws.Range("L99").Value = 10896.85
$ws.Range("M99").Value = -2042
$ws.Range("N99").Value = -13892.85

$ws.Range("H126").Value = 9425.48
$ws.Range("I126").Value = 3540
$ws.Range("J126").Value = 10896.85
$ws.Range("K126").Value = 10620
$ws.Range("L126").Value = 32690.55
$ws.Range("M126").Value = -8150
$ws.Range("N126").Value = -37630.55

$ws.Range("H136").Value = 18773.84
$ws.Range("I136").Value = 8648.154
$ws.Range("K136").Value = 25944.462
$ws.Range("M136").Value = -23394.462

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 79164360
$ws.Range("I4").Value = 114126280
$ws.Range("J4").Value = 500048
$ws.Range("K4").Value = 342378840
$ws.Range("L4").Value = 1500144
$ws.Range("M4").Value = -342378728
$ws.Range("N4").Value = -1500368

$ws.Range("H11").Value = 1070.5385
$ws.Range("I11").Value = 1231.4445
$ws.Range("K11").Value = 3694.3335
$ws.Range("M11").Value = -3554.3335

$ws.Range("H12").Value = 83.42856999999999
$ws.Range("I12").Value = 110.333336
$ws.Range("J12").Value = 63.25
$ws.Range("K12").Value = 331.000008
$ws.Range("L12").Value = 189.75
$ws.Range("M12").Value = -158.000008
$ws.Range("N12").Value = -535.75

$ws.Range("H33").Value = 442.33334
$ws.Range("I33").Value = 157.33333
$ws.Range("K33").Value = 943.9999799999999
$ws.Range("M33").Value = -660.9999799999999

$ws.Range("H44").Value = 750
$ws.Range("J44").Value = 1000
$ws.Range("L44").Value = 3000
$ws.Range("N44").Value = -3796

$ws.Range("H68").Value = 6042.077
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 6042.077
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 18126.231
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -19748.231

$ws.Range("H71").Value = 6042.077
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 6042.077
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 54378.693
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -62490.693

$ws.Range("H113").Value = 1534.9048
$ws.Range("I113").Value = 1194.6666
$ws.Range("J113").Value = 1671
$ws.Range("K113").Value = 3583.9998
$ws.Range("L113").Value = 5013
$ws.Range("M113").Value = -1413.9998
$ws.Range("N113").Value = -9353

$ws.Range("H132").Value = 2619.8
$ws.Range("J132").Value = 2899.5
$ws.Range("L132").Value = 26095.5
$ws.Range("N132").Value = -31155.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 41230.76
$ws.Range("J108").Value = 41230.76
$ws.Range("L108").Value = 41230.76
$ws.Range("N108").Value = -48910.76

$ws.Range("H137").Value = 55999.668
$ws.Range("J137").Value = 55999.668
$ws.Range("L137").Value = 55999.668
$ws.Range("N137").Value = -66199.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6676.7334
$ws.Range("I22").Value = 1770.7142
$ws.Range("K22").Value = 1770.7142
$ws.Range("M22").Value = -1475.7142

$ws.Range("H27").Value = 6676.7334
$ws.Range("I27").Value = 1770.7142
$ws.Range("K27").Value = 1770.7142
$ws.Range("M27").Value = -1663.7142

$ws.Range("H93").Value = 11836.083
$ws.Range("I93").Value = 10604.866
$ws.Range("J93").Value = 13888.111
$ws.Range("K93").Value = 10604.866
$ws.Range("L93").Value = 13888.111
$ws.Range("M93").Value = -9356.866
$ws.Range("N93").Value = -16384.111

$ws.Range("H94").Value = 62000
$ws.Range("J94").Value = 62000
$ws.Range("L94").Value = 62000
$ws.Range("N94").Value = -63352

$ws.Range("H136").Value = 12597.325
$ws.Range("I136").Value = 13038.579
$ws.Range("J136").Value = 12248
$ws.Range("K136").Value = 39115.737
$ws.Range("L136").Value = 36744
$ws.Range("M136").Value = -36565.737
$ws.Range("N136").Value = -41844

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 9697.419
$ws.Range("I132").Value = 3842
$ws.Range("J132").Value = 16807.572
$ws.Range("K132").Value = 11526
$ws.Range("L132").Value = 50422.716
$ws.Range("M132").Value = -8996
$ws.Range("N132").Value = -55482.716

$ws.Range("H133").Value = 76571
$ws.Range("J133").Value = 76571
$ws.Range("L133").Value = 76571
$ws.Range("N133").Value = -86691
